$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the game 2021020003 data (previously row 4), with an updated
# Teams description that includes each team's record.
# Force the GameID to be stored as text (it's an identifier, not a number)
# the same way the original workbook stores it, instead of letting Excel's
# automatic type inference turn the numeric-looking string into a number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2021020003"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "Scotiabank Arena"
$ws.Range("C2").Value = "86.0 meters meters"
$ws.Range("D2").Value = "October 13, 2021"
$ws.Range("E2").Value = "Montreal Canadiens (0-0-0) vs. Toronto Maple Leafs (0-0-0)"
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = "00:03"
$ws.Range("H2").Value = "SHOT"
$ws.Range("I2").Value = "Jeff Petry (Canadiens)"
$ws.Range("J2").Value = "23.41 feet"
$ws.Range("K2").Value = "1 - 2"
$ws.Range("L2").Value = "https://www.nhl.com/scores/htmlreports/20212022/PL020003.HTM"

# Remove the old rows 3 and 4 (game 2021020001 second shot, and the
# duplicate game 2021020003 row), shrinking the sheet down to A1:L2.
$ws.Range("A4:L4").EntireRow.Delete()
$ws.Range("A3:L3").EntireRow.Delete()
